$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet, add the new NOTES sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Problem-soln"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "NOTES"

# --- 2. Problem-soln: bold the "Hashing Basic" topic cell (A3) ---
$ws1.Range("A3").Font.Bold = $true

# --- 3. Problem-soln: append a new problem row (row 5) ---
$ws1.Range("A3").Copy()
$ws1.Range("A5").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("A5").Font.Bold = $false
$ws1.Range("A5").Value = "Count min max frequency "

$ws1.Range("B4").Copy()
$ws1.Range("B5").PasteSpecial(-4122)   # xlPasteFormats (hyperlink style)
$ws1.Range("B5").Value = "https://www.codingninjas.com/studio/problems/k-most-occurrent-numbers_625382?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf"
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://www.codingninjas.com/studio/problems/k-most-occurrent-numbers_625382?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf") | Out-Null

$ws1.Range("C4").Copy()
$ws1.Range("C5").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("C5").Value = "hash as map. Set returun array vals as -1. Then loop over hashmap to get min and max frequency. In first itr set  array min max to first element . Bit tricky but easy."

$ws1.Rows.Item(5).RowHeight = 75

$ws1.Range("C5").Select()

# --- 4. NOTES sheet: headers (row 1) ---
$ws2.Range("A1").Value = "TOPIC"
$ws2.Range("B1").Value = "Details"

$ws1.Range("A3").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Rows.Item(1).Font.Bold = $true

# --- 5. NOTES sheet: first note entry (row 3) ---
$ws2.Range("A3").Value = "Hashing "
$ws2.Range("B3").Value = "Data can be hashed in the form of array /map to store their frequency or other data. For large Data hasing technique such as division hasing can be used. Like mod element with 10 and store in map of <int,LinkedList>  e.g 8-> 8,18,28,38..."

$ws1.Range("A4").Copy()
$ws2.Range("B3").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Rows.Item(3).RowHeight = 45

# --- 6. NOTES sheet: column widths ---
$ws2.Columns.Item(1).ColumnWidth = 55.45182291666667
$ws2.Columns.Item(2).ColumnWidth = 79.45182291666667

$ws2.Range("B3").Select()
